$d = $word.ActiveDocument

# This document's "Comments" Heading-1 paragraph currently carries Word's
# cached <w:lastRenderedPageBreak/> rendering marker on its run, while the
# following "Unresolved Issue" Heading-1 paragraph does not. The edit moves
# that marker: it is cleared from the "Comments" heading run and placed on
# the "Unresolved Issue" heading run instead. No visible text changes.

$commentsHeadingIndex = 0
$unresolvedHeadingIndex = 0

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs($i)
    $styleName = $para.Style.NameLocal
    $text = $para.Range.Text.TrimEnd("`r", "`a")

    if ($styleName -eq "Heading 1" -and $text -eq "Comments" -and $commentsHeadingIndex -eq 0) {
        $commentsHeadingIndex = $i
    }
    if ($styleName -eq "Heading 1" -and $text -eq "Unresolved Issue" -and $unresolvedHeadingIndex -eq 0) {
        $unresolvedHeadingIndex = $i
    }
}

# Step 1: rewrite the "Comments" heading run's text in place (same text);
# this regenerates the run without the stale <w:lastRenderedPageBreak/>
# marker while leaving the paragraph's style/formatting untouched.
if ($commentsHeadingIndex -gt 0) {
    $commentsRange = $d.Paragraphs($commentsHeadingIndex).Range
    $commentsRange.Text = "Comments"
}

# Step 2: rebuild the "Unresolved Issue" heading paragraph so the run
# carries <w:lastRenderedPageBreak/> immediately before its text, matching
# the paragraph style / run language it already had.
if ($unresolvedHeadingIndex -gt 0) {
    $issuePara = $d.Paragraphs($unresolvedHeadingIndex)
    $issueText = $issuePara.Range.Text.TrimEnd("`r", "`a")
    $issueStyleId = $issuePara.Style.NameLocal -replace ' ', ''
    $issueLang = $issuePara.Range.LanguageID
    if (-not $issueLang) { $issueLang = "en-US" }

    $issueXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="' + $issueStyleId + '"/><w:rPr><w:lang w:val="' + $issueLang + '"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="' + $issueLang + '"/></w:rPr><w:lastRenderedPageBreak/><w:t>' + $issueText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $issuePara.Range.InsertXML($issueXml)
}
